$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.894.49'

$ws.Range('E2').Value = '  +7.18%  '

$ws.Range('D3').Value = '3.859.57'

$ws.Range('E3').Value = '  +13.51%  '

$ws.Range('E4').Value = '  +0.21%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '426.27'
$ws.Range('D5').Style = 'Normal'

$ws.Range('E5').Value = '  +11.87%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '130.59'
$ws.Range('D6').Style = 'Normal'

$ws.Range('E6').Value = '  +8.48%  '

$ws.Range('D7').Value = '3.852.03'

$ws.Range('E7').Value = '  +9.25%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.611'
$ws.Range('D8').Style = 'Normal'

$ws.Range('E8').Value = '  +6.94%  '

$ws.Range('E9').Value = '  -0.08%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.725'
$ws.Range('D10').Style = 'Normal'

$ws.Range('E10').Value = '  +11.99%  '

$ws.Range('E11').Value = '  +16.27%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000346'
$ws.Range('D12').Style = 'Normal'

$ws.Range('E12').Value = '  +23.38%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '40.99'
$ws.Range('D13').Style = 'Normal'

$ws.Range('E13').Value = '  +8.78%  '

$ws.Range('D14').Value = '4.478.07'

$ws.Range('E14').Value = '  +14.72%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '10.24'
$ws.Range('D15').Style = 'Normal'

$ws.Range('E15').Value = '  +14.88%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '15.95'
$ws.Range('D16').Style = 'Normal'

$ws.Range('E16').Value = '  +33.17%  '

$ws.Range('D17').Value = '3.892.43'

$ws.Range('E17').Value = '  +14.98%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.137'
$ws.Range('D18').Style = 'Normal'

$ws.Range('E18').Value = '  +1.32%  '

$ws.Range('E19').Value = '  +10.09%  '

$ws.Range('D20').Value = '67.030.52'

$ws.Range('E20').Value = '  +7.96%  '

$ws.Range('E21').Value = '  +9.85%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '412.33'
$ws.Range('D22').Style = 'Normal'

$ws.Range('E22').Value = '  +8.24%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '14.90'
$ws.Range('D23').Style = 'Normal'

$ws.Range('E23').Value = '  +13.64%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '84.10'
$ws.Range('D24').Style = 'Normal'

$ws.Range('E24').Value = '  +7.13%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.02'
$ws.Range('D25').Style = 'Normal'

$ws.Range('E25').Value = '  +11.17%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '37.77'
$ws.Range('D26').Style = 'Normal'

$ws.Range('E26').Value = '  +17.89%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.99'
$ws.Range('D27').Style = 'Normal'

$ws.Range('E27').Value = '  +17.22%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '3.24'
$ws.Range('D28').Style = 'Normal'

$ws.Range('E28').Value = '  +12.26%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.34'
$ws.Range('D29').Style = 'Normal'

$ws.Range('E29').Value = '  +3.34%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.88'
$ws.Range('D30').Style = 'Normal'

$ws.Range('E30').Value = '  +36.68%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '724.66'
$ws.Range('D31').Style = 'Normal'

$ws.Range('E31').Value = '  +13.95%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '13.62'
$ws.Range('D32').Style = 'Normal'

$ws.Range('E32').Value = '  +19.01%  '

$ws.Range('E33').Value = '  +14.91%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.78'
$ws.Range('D34').Style = 'Normal'

$ws.Range('E34').Value = '  +8.53%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.997'
$ws.Range('D35').Style = 'Normal'

$ws.Range('E35').Value = '  -0.13%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '38.99'
$ws.Range('D36').Style = 'Normal'

$ws.Range('E36').Value = '  +10.27%  '

$ws.Range('E37').Value = '  +4.79%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '55.99'
$ws.Range('D38').Style = 'Normal'

$ws.Range('B39').Value = 'PEPE'

$ws.Range('C39').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'

$ws.Range('D39').Value = '0.0₃0764'

$ws.Range('E39').Value = '  +31.90%  '

$ws.Range('B40').Value = 'NEARProtocol'

$ws.Range('C40').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.26'
$ws.Range('D40').Style = 'Normal'

$ws.Range('E40').Value = '  +34.67%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0460'
$ws.Range('D41').Style = 'Normal'

$ws.Range('E41').Value = '  +9.85%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.87'
$ws.Range('D42').Style = 'Normal'

$ws.Range('E42').Value = '  +12.77%  '

$ws.Range('E43').Value = '  +1.45%  '

$ws.Range('B44').Value = 'Stellar'

$ws.Range('C44').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.135'
$ws.Range('D44').Style = 'Normal'

$ws.Range('E44').Value = '  +5.44%  '

$ws.Range('B45').Value = 'ApeXProtocol'

$ws.Range('C45').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.20'
$ws.Range('D45').Style = 'Normal'

$ws.Range('E45').Value = '  +9.92%  '

$ws.Range('E46').Value = '  +13.25%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.312'
$ws.Range('D47').Style = 'Normal'

$ws.Range('E47').Value = '  +16.60%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '142.19'
$ws.Range('D48').Style = 'Normal'

$ws.Range('E48').Value = '  +5.04%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.03'
$ws.Range('D49').Style = 'Normal'

$ws.Range('E49').Value = '  +8.09%  '

$ws.Range('E50').Value = '  +9.86%  '

$ws.Range('E51').Value = '  +6.62%  '
